$d = $word.ActiveDocument

# --- 1. Fix the date: 2024-10-02 -> 2024-10-09 ---
$d.Content.Find.Execute("2024-10-02", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-10-09", 2)

# --- 2. Fill in the empty log entries in the second table ---
$t1 = $d.Tables.Item(2)

# "What was accomplished last week?" section
$t1.Cell(4, 1).Range.Text  = "Converted map pdf to images."
$t1.Cell(5, 1).Range.Text  = "Created functioning code in Godot."
$t1.Cell(6, 1).Range.Text  = "Continued working on Functional Specification."

# "What will be done next week?" section
$t1.Cell(15, 1).Range.Text = "Continue working on Functional Specification."
$t1.Cell(16, 1).Range.Text = "Continue working on Research Poster."
